$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) In the "Metric for evaluating" bullet list, insert a new sub-bullet
#    "LOOK AT COST FUNCTION" right before "Jaccard" (i.e. right after "Dice").
#    This pushes "Jaccard" / "Modes of failure" / "Speed of convergence" down
#    by one position, matching the text-shift pattern in the diff.
# ---------------------------------------------------------------------------
$diceRun = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq "Dice") {
        $diceRun = $p
    }
}

$diceRun.Range.InsertParagraphAfter()
$newBullet = $diceRun.Next()
$newBullet.Range.Text = "LOOK AT COST FUNCTION"

# ---------------------------------------------------------------------------
# 2) After "What are we trying to segment?" insert the new block of notes,
#    before the trailing bookmark paragraph:
#      Simulated image and brain image
#      <blank>
#      Compare SD and GA for easy and hard questions   (two runs)
#      <blank>
#      What kind of segmentation formulation will you use?
#      <blank>
# ---------------------------------------------------------------------------
$segmentPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq "What are we trying to segment?") {
        $segmentPara = $p
    }
}

# "Simulated image and brain image"
$segmentPara.Range.InsertParagraphAfter()
$pSimulated = $segmentPara.Next()
$pSimulated.Range.Text = "Simulated image and brain image"

# blank separator paragraph
$pSimulated.Range.InsertParagraphAfter()
$pBlank1 = $pSimulated.Next()

# "Compare SD and GA " + "for easy and hard questions" (kept as two runs)
$pBlank1.Range.InsertParagraphAfter()
$pCompareA = $pBlank1.Next()
$pCompareA.Range.Text = "Compare SD and GA "

$pCompareA.Range.InsertParagraphAfter()
$pCompareB = $pCompareA.Next()
$pCompareB.Range.Text = "for easy and hard questions"

# merge pCompareA and pCompareB into a single paragraph with two runs by
# deleting the paragraph mark between them
$mark = $d.Range($pCompareA.Range.End - 1, $pCompareA.Range.End)
$mark.Delete()

# blank separator paragraph
$pCompareA.Range.InsertParagraphAfter()
$pBlank2 = $pCompareA.Next()

# "What kind of segmentation formulation will you use?"
$pBlank2.Range.InsertParagraphAfter()
$pFormulation = $pBlank2.Next()
$pFormulation.Range.Text = "What kind of segmentation formulation will you use?"

# trailing blank separator paragraph
$pFormulation.Range.InsertParagraphAfter()
$pBlank3 = $pFormulation.Next()
